$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-16 (topics -1 through 13) with the reduced candidate-topic data
$data = @(
    @(-1, 258, '-1_ülke_genç_iş_kardeş', 'ülke, genç, iş, kardeş, millet'),
    @(0, 162, '0_buluş_millet_muhteşem_erzurum', 'buluş, millet, muhteşem, erzurum, hatay'),
    @(1, 106, '1_konut_deprem_temel_at', 'konut, deprem, temel, at, ev'),
    @(2, 79, '2_rahmet_atatürk_şehit_dile', 'rahmet, atatürk, şehit, dile, allah'),
    @(3, 71, '3_türkiye_yüzyıl_ülke_mayıs', 'türkiye, yüzyıl, ülke, mayıs, millet'),
    @(4, 62, '4_seçim_sandık_millet_mayıs', 'seçim, sandık, millet, mayıs, oy'),
    @(5, 45, '5_teşekkür_kardeş_muhteşem_güzel', 'teşekkür, kardeş, muhteşem, güzel, istanbul'),
    @(6, 43, '6_canlı_yayın_tv_basın', 'canlı, yayın, tv, basın, açıkla'),
    @(7, 36, '7_gel_sev_haydi_vatan', 'gel, sev, haydi, vatan, bura'),
    @(8, 30, '8_milyar_yatırım_lira_kamu', 'milyar, yatırım, lira, kamu, dolar'),
    @(9, 24, '9_nükleer_santral_enerji_üretim', 'nükleer, santral, enerji, üretim, ülke'),
    @(10, 15, '10_emekli_maaş_hanım_aile', 'emekli, maaş, hanım, aile, düşük'),
    @(11, 12, '11_uzay_altay_milli_tank', 'uzay, altay, milli, tank, uçak'),
    @(12, 12, '12_öğretmen_ata_eğitim_engelli', 'öğretmen, ata, eğitim, engelli, okul'),
    @(13, 10, '13_para_yatırım_erdoğan_çalış', 'para, yatırım, erdoğan, çalış, ver'),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Remove the now-unused rows (17-25), shifting dimension down to A1:D16
$ws.Rows("17:25").Delete()
